$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '61.115.97'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  -1.63%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.970.95'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  -1.49%  '
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  +0.08%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '593.49'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +1.78%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '142.58'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  -2.96%  '
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  +0.01%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '2.972.20'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  -1.41%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.512'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  -2.35%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.147'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  -1.44%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '6.00'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +3.06%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.451'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  +1.71%  '
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  -1.09%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '33.93'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  -1.72%  '
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  +2.95%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.464.04'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  -1.25%  '
$ws.Range("B17").Value = 'Polkadot'
$ws.Range("C17").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '6.87'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  -2.27%  '
$ws.Range("B18").Value = 'WrappedBTC'
$ws.Range("C18").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '61.135.22'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  -1.49%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '2.972.44'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  -1.23%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '445.79'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  -3.94%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '13.93'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  +0.04%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.679'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  -0.77%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.29'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  -2.30%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '81.33'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  -0.57%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '10.55'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  +6.34%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.17'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  -3.45%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '11.92'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  -2.69%  '
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  -0.02%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.67'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  +1.92%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.999'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  +0.02%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '7.14'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -0.28%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '2.04'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  -2.03%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '27.01'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  -5.71%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.110'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  +0.71%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.0₃0805'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  +1.48%  '
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  -1.23%  '
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  -0.56%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '49.99'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  -0.54%  '
$ws.Range("B39").Value = 'Cosmos'
$ws.Range("C39").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '8.93'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  -1.14%  '
$ws.Range("B40").Value = 'Stacks'
$ws.Range("C40").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.01'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -4.43%  '
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +7.01%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.82'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  -3.83%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '382.76'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  +0.06%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0348'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  -1.56%  '
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  -2.52%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '38.24'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  +1.60%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.676.75'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  -2.85%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '130.17'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  +1.75%  '
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  +0.15%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.107'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  -1.50%  '
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  -1.46%  '
